# Applies the "Note some obvious coding mistakes" revision:
#   1. Sumstat -> Sum + (ins "maryS") + (del "s") + tat   (tracked change, author "Miklos Koren")
#   2. Removes the stray "_GoBack" bookmark around "Figure I, Table..."
#   3. Registers the (otherwise-unused-until-now) built-in "Revision" paragraph style,
#      which Word silently adds to styles.xml once track-changes edits exist in the doc.

$d = $word.ActiveDocument

# Track-change edits must be attributed to the same author as in the target revision.
$word.UserName = "Miklos Koren"

# --- 1. "Sumstat" -> "Sum" + ins("maryS") + del("s") + "tat" -------------------------
$r = $d.Content
$found = $r.Find.Execute("Sumstat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $wordStart = $r.Start

    $d.TrackRevisions = $true

    # Insert "maryS" right after "Sum" (3 chars in), before the "s" of "Sumstat".
    $insPoint = $d.Range($wordStart + 3, $wordStart + 3)
    $insPoint.InsertBefore("maryS")

    # Delete the old "s" that used to directly follow "Sum" (now shifted right by
    # the 5 characters we just inserted).
    $delRange = $d.Range($wordStart + 3 + 5, $wordStart + 3 + 5 + 1)
    $delRange.Delete()

    $d.TrackRevisions = $false
}

# --- 2. Remove the leftover "_GoBack" bookmark ----------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Register the "Revision" paragraph style (added by Word once tracked changes
#        with paragraph-mark deletions exist in the document) -------------------------
$revStyle = $d.Styles.Add("Revision", 1)
$revStyle.Priority = 99
$revStyle.Visibility = $true
$revStyle.ParagraphFormat.SpaceAfter = 0
$revStyle.ParagraphFormat.LineSpacingRule = 0

Write-Host "done"
